$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.598.45"
$ws.Range("E2").Value = "  +5.93%  "
$ws.Range("D3").Value = "2.298.38"
$ws.Range("E3").Value = "  +3.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'304.52"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "'100.70"
$ws.Range("E6").Value = "  +11.39%  "
$ws.Range("D7").Value = "'0.566"
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = "  +6.63%  "
$ws.Range("D10").Value = "'36.57"
$ws.Range("E10").Value = "  +10.54%  "
$ws.Range("D11").Value = "'0.0792"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("E12").Value = "  +6.27%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "2.649.85"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").Value = "2.300.37"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").Value = "'13.85"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("E17").Value = "  +4.75%  "
$ws.Range("D18").Value = "46.608.49"
$ws.Range("E18").Value = "  +6.34%  "
$ws.Range("D19").Value = "'13.06"
$ws.Range("E19").Value = "  +10.78%  "
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").Value = "  +3.82%  "
$ws.Range("D21").Value = "'6.03"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").Value = "'66.27"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("D23").Value = "'248.83"
$ws.Range("E23").Value = "  +5.35%  "
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  +3.99%  "
$ws.Range("D27").Value = "'42.74"
$ws.Range("E27").Value = "  +9.85%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'9.85"
$ws.Range("E29").Value = "  +5.18%  "
$ws.Range("D30").Value = "'20.00"
$ws.Range("E30").Value = "  +4.35%  "
$ws.Range("E31").Value = "  +12.82%  "
$ws.Range("D32").Value = "'5.66"
$ws.Range("E32").Value = "  +4.65%  "
$ws.Range("D33").Value = "'147.35"
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("E34").Value = "  +4.46%  "
$ws.Range("E35").Value = "  +15.78%  "
$ws.Range("E36").Value = "  +11.25%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  +6.07%  "
$ws.Range("D39").Value = "'16.08"
$ws.Range("E39").Value = "  +20.75%  "
$ws.Range("D40").Value = "'4.00"
$ws.Range("E40").Value = "  +10.78%  "
$ws.Range("D41").Value = "'3.36"
$ws.Range("E41").Value = "  +6.32%  "
$ws.Range("D42").Value = "'0.0302"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  +10.80%  "
$ws.Range("D45").Value = "1.820.14"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").Value = "'88.18"
$ws.Range("E46").Value = "  +21.02%  "
$ws.Range("E47").Value = "  +5.97%  "
$ws.Range("D48").Value = "'73.38"
$ws.Range("E48").Value = "  +8.37%  "
$ws.Range("D49").Value = "'4.91"
$ws.Range("E49").Value = "  +7.17%  "
$ws.Range("D50").Value = "'96.09"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("D51").Value = "'53.76"
$ws.Range("E51").Value = "  +5.35%  "
